$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9211618304252625
$ws.Range("B1").Value = 2.683574676513672
$ws.Range("C1").Value = 2.539758682250977
$ws.Range("D1").Value = 2.541419267654419
$ws.Range("E1").Value = 1.889025568962097
